# Update 'F' column (想去人数 / interest counts) values across sheets
# per gh-pages regeneration at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 599
$ws.Range("F3").Value = 5700
$ws.Range("F4").Value = 66
$ws.Range("F5").Value = 460
$ws.Range("F7").Value = 1013
$ws.Range("F8").Value = 385
$ws.Range("F12").Value = 3100
$ws.Range("F13").Value = 1930
$ws.Range("F15").Value = 60
$ws.Range("F17").Value = 57
$ws.Range("F18").Value = 161
$ws.Range("F22").Value = 53
$ws.Range("F23").Value = 26
$ws.Range("F24").Value = 3595
$ws.Range("F25").Value = 1145
$ws.Range("F26").Value = 2857
$ws.Range("F27").Value = 286
$ws.Range("F28").Value = 2222
$ws.Range("F29").Value = 4154
$ws.Range("F30").Value = 110
$ws.Range("F31").Value = 923
$ws.Range("F32").Value = 473
$ws.Range("F33").Value = 1321
$ws.Range("F34").Value = 77
$ws.Range("F35").Value = 31
$ws.Range("F36").Value = 1015
$ws.Range("F37").Value = 1284
$ws.Range("F38").Value = 67
$ws.Range("F39").Value = 1076
$ws.Range("F40").Value = 695
$ws.Range("F41").Value = 572
$ws.Range("F42").Value = 420
$ws.Range("F43").Value = 7
$ws.Range("F44").Value = 73
$ws.Range("F45").Value = 320
$ws.Range("F46").Value = 3584

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 22
$ws.Range("F10").Value = 909
$ws.Range("F16").Value = 15
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = 3
$ws.Range("F22").Value = 2
$ws.Range("F23").Value = 39
$ws.Range("F25").Value = 14

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 599
$ws.Range("F3").Value = 5700
$ws.Range("F4").Value = 66
$ws.Range("F5").Value = 22
$ws.Range("F8").Value = 385
$ws.Range("F10").Value = 3100
$ws.Range("F12").Value = 1930
$ws.Range("F14").Value = 60
$ws.Range("F16").Value = 909
$ws.Range("F18").Value = 161
$ws.Range("F21").Value = 3595
$ws.Range("F22").Value = 15
$ws.Range("F24").Value = 1145
$ws.Range("F25").Value = 3
$ws.Range("F26").Value = 2857
$ws.Range("F27").Value = 2222
$ws.Range("F28").Value = 4154
$ws.Range("F29").Value = 3
$ws.Range("F30").Value = 110
$ws.Range("F31").Value = 923
$ws.Range("F32").Value = 1321
$ws.Range("F33").Value = 31
$ws.Range("F34").Value = 1015
$ws.Range("F35").Value = 1284
$ws.Range("F36").Value = 67
$ws.Range("F37").Value = 1076
$ws.Range("F39").Value = 695
$ws.Range("F40").Value = 2
$ws.Range("F41").Value = 420
$ws.Range("F42").Value = 39
$ws.Range("F44").Value = 14
$ws.Range("F45").Value = 73
$ws.Range("F47").Value = 320
$ws.Range("F48").Value = 3584
